$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 105012296
$ws.Range("I86").Value = 131264570
$ws.Range("J86").Value = 3199.75
$ws.Range("K86").Value = 131264570
$ws.Range("L86").Value = 3199.75
$ws.Range("M86").Value = -131263447
$ws.Range("N86").Value = -5445.75
$ws.Range("H89").Value = 105012296
$ws.Range("I89").Value = 131264570
$ws.Range("J89").Value = 3199.75
$ws.Range("K89").Value = 656322850
$ws.Range("L89").Value = 15998.75
$ws.Range("M89").Value = -656317234
$ws.Range("N89").Value = -27230.75
$ws.Range("H92").Value = 22224066
$ws.Range("I92").Value = 27779250
$ws.Range("J92").Value = 3333.3333
$ws.Range("K92").Value = 27779250
$ws.Range("L92").Value = 3333.3333
$ws.Range("M92").Value = -27778002
$ws.Range("N92").Value = -5829.3333
$ws.Range("H112").Value = 6781.6816
$ws.Range("I112").Value = 970
$ws.Range("J112").Value = 7699.316
$ws.Range("K112").Value = 2910
$ws.Range("L112").Value = 23097.948
$ws.Range("M112").Value = -1802
$ws.Range("N112").Value = -25313.948
$ws.Range("H137").Value = 4764.4287
$ws.Range("J137").Value = 2460
$ws.Range("L137").Value = 7380
$ws.Range("N137").Value = -12480
$ws.Range("H138").Value = 193352.73
$ws.Range("I138").Value = 8611.75
$ws.Range("J138").Value = 207295.45
$ws.Range("K138").Value = 25835.25
$ws.Range("L138").Value = 621886.3500000001
$ws.Range("M138").Value = -20695.25
$ws.Range("N138").Value = -632166.3500000001
$ws.Range("H141").Value = 3204.532
$ws.Range("I141").Value = 1716.7742
$ws.Range("J141").Value = 6087.0625
$ws.Range("K141").Value = 5150.3226
$ws.Range("L141").Value = 18261.1875
$ws.Range("M141").Value = 29.67739999999958
$ws.Range("N141").Value = -28621.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2903.7896
$ws.Range("I61").Value = 1997.6666
$ws.Range("J61").Value = 4457.143
$ws.Range("K61").Value = 1997.6666
$ws.Range("L61").Value = 4457.143
$ws.Range("M61").Value = -1785.6666
$ws.Range("N61").Value = -4881.143
$ws.Range("H74").Value = 1036
$ws.Range("I74").Value = 763.6957
$ws.Range("J74").Value = 1605.3636
$ws.Range("K74").Value = 763.6957
$ws.Range("L74").Value = 1605.3636
$ws.Range("M74").Value = 110.3043
$ws.Range("N74").Value = -3353.3636
$ws.Range("H76").Value = 54666.668
$ws.Range("J76").Value = 54666.668
$ws.Range("L76").Value = 54666.668
$ws.Range("N76").Value = -55342.668
$ws.Range("H77").Value = 1036
$ws.Range("I77").Value = 763.6957
$ws.Range("J77").Value = 1605.3636
$ws.Range("K77").Value = 3818.4785
$ws.Range("L77").Value = 8026.817999999999
$ws.Range("M77").Value = 549.5214999999998
$ws.Range("N77").Value = -16762.818
$ws.Range("H79").Value = 54666.668
$ws.Range("J79").Value = 54666.668
$ws.Range("L79").Value = 54666.668
$ws.Range("N79").Value = -57006.668
$ws.Range("H102").Value = 1010
$ws.Range("I102").Value = 1010
$ws.Range("K102").Value = 1010
$ws.Range("M102").Value = 612
$ws.Range("H110").Value = 1719.7333
$ws.Range("I110").Value = 1615.0769
$ws.Range("K110").Value = 1615.0769
$ws.Range("M110").Value = 429.9231
$ws.Range("H136").Value = 2903.7896
$ws.Range("I136").Value = 1997.6666
$ws.Range("J136").Value = 4457.143
$ws.Range("K136").Value = 5992.9998
$ws.Range("L136").Value = 13371.429
$ws.Range("M136").Value = -3442.9998
$ws.Range("N136").Value = -18471.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66668980
$ws.Range("I86").Value = 90910840
$ws.Range("J86").Value = 3875
$ws.Range("K86").Value = 90910840
$ws.Range("L86").Value = 3875
$ws.Range("M86").Value = -90909717
$ws.Range("N86").Value = -6121
$ws.Range("H89").Value = 66668980
$ws.Range("I89").Value = 90910840
$ws.Range("J89").Value = 3875
$ws.Range("K89").Value = 454554200
$ws.Range("L89").Value = 19375
$ws.Range("M89").Value = -454548584
$ws.Range("N89").Value = -30607
$ws.Range("H94").Value = 492.54166
$ws.Range("I94").Value = 430.3889
$ws.Range("J94").Value = 679
$ws.Range("K94").Value = 430.3889
$ws.Range("L94").Value = 679
$ws.Range("M94").Value = 20.61110000000002
$ws.Range("N94").Value = -1581
$ws.Range("H134").Value = 3539.4583
$ws.Range("I134").Value = 3116.6428
$ws.Range("J134").Value = 4131.4
$ws.Range("K134").Value = 9349.928400000001
$ws.Range("L134").Value = 12394.2
$ws.Range("M134").Value = -6814.928400000001
$ws.Range("N134").Value = -17464.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4517.0435
$ws.Range("I31").Value = 1062.9286
$ws.Range("K31").Value = 1062.9286
$ws.Range("M31").Value = -767.9286
$ws.Range("H34").Value = 4517.0435
$ws.Range("I34").Value = 1062.9286
$ws.Range("K34").Value = 1062.9286
$ws.Range("M34").Value = -860.9286
$ws.Range("H58").Value = 1502.3889
$ws.Range("I58").Value = 1246.4286
$ws.Range("J58").Value = 1665.2727
$ws.Range("K58").Value = 1246.4286
$ws.Range("L58").Value = 1665.2727
$ws.Range("M58").Value = -1043.4286
$ws.Range("N58").Value = -2071.2727
$ws.Range("H132").Value = 9261275
$ws.Range("I132").Value = 1067.8
$ws.Range("J132").Value = 20836534
$ws.Range("K132").Value = 3203.4
$ws.Range("L132").Value = 62509602
$ws.Range("M132").Value = -673.3999999999996
$ws.Range("N132").Value = -62514662
$ws.Range("H134").Value = 3295.4119
$ws.Range("I134").Value = 2167
$ws.Range("J134").Value = 4564.875
$ws.Range("K134").Value = 6501
$ws.Range("L134").Value = 13694.625
$ws.Range("M134").Value = -3966
$ws.Range("N134").Value = -18764.625
$ws.Range("H136").Value = 1502.3889
$ws.Range("I136").Value = 1246.4286
$ws.Range("J136").Value = 1665.2727
$ws.Range("K136").Value = 3739.2858
$ws.Range("L136").Value = 4995.8181
$ws.Range("M136").Value = -1189.2858
$ws.Range("N136").Value = -10095.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 982.96
$ws.Range("I5").Value = 554.95
$ws.Range("K5").Value = 1664.85
$ws.Range("M5").Value = -1552.85
$ws.Range("H33").Value = 8499.916999999999
$ws.Range("I33").Value = 12687.375
$ws.Range("J33").Value = 125
$ws.Range("K33").Value = 76124.25
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -75841.25
$ws.Range("N33").Value = -1316
$ws.Range("H68").Value = 3902
$ws.Range("I68").Value = 3902
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 11706
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -10895
$ws.Range("H71").Value = 3902
$ws.Range("I71").Value = 3902
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 35118
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -31062
$ws.Range("H113").Value = 968.7179599999999
$ws.Range("I113").Value = 705.55554
$ws.Range("J113").Value = 1560.8334
$ws.Range("K113").Value = 2116.66662
$ws.Range("L113").Value = 4682.5002
$ws.Range("M113").Value = 53.33338000000003
$ws.Range("N113").Value = -9022.5002
$ws.Range("H131").Value = 954.04877
$ws.Range("J131").Value = 1023.9459
$ws.Range("L131").Value = 3071.8377
$ws.Range("N131").Value = -13151.8377
$ws.Range("H135").Value = 982.96
$ws.Range("I135").Value = 554.95
$ws.Range("K135").Value = 4994.55
$ws.Range("M135").Value = -2459.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2449.923
$ws.Range("I102").Value = 2294.9
$ws.Range("J102").Value = 2966.6667
$ws.Range("K102").Value = 2294.9
$ws.Range("L102").Value = 2966.6667
$ws.Range("M102").Value = -672.9000000000001
$ws.Range("N102").Value = -6210.6667
$ws.Range("H122").Value = 4404.3335
$ws.Range("I122").Value = 3010.3333
$ws.Range("K122").Value = 9030.999899999999
$ws.Range("M122").Value = -6580.999899999999
$ws.Range("H126").Value = 2852.4
$ws.Range("I126").Value = 2815.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8446.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5976.5
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 76926960
$ws.Range("I7").Value = 83336960
$ws.Range("K7").Value = 83336960
$ws.Range("M7").Value = -83336848
$ws.Range("H82").Value = 2291.3635
$ws.Range("J82").Value = 2133.6667
$ws.Range("L82").Value = 2133.6667
$ws.Range("N82").Value = -2855.6667
$ws.Range("H85").Value = 2291.3635
$ws.Range("J85").Value = 2133.6667
$ws.Range("L85").Value = 2133.6667
$ws.Range("N85").Value = -4629.6667
$ws.Range("H126").Value = 76926960
$ws.Range("I126").Value = 83336960
$ws.Range("K126").Value = 250010880
$ws.Range("M126").Value = -250008410

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 47261.832
$ws.Range("J69").Value = 47261.832
$ws.Range("L69").Value = 47261.832
$ws.Range("N69").Value = -48759.832
$ws.Range("H72").Value = 47261.832
$ws.Range("J72").Value = 47261.832
$ws.Range("L72").Value = 141785.496
$ws.Range("N72").Value = -149273.496
$ws.Range("H132").Value = 4632432.5
$ws.Range("I132").Value = 3195.5264
$ws.Range("J132").Value = 9806285
$ws.Range("K132").Value = 9586.5792
$ws.Range("L132").Value = 29418855
$ws.Range("M132").Value = -7056.5792
$ws.Range("N132").Value = -29423915
$ws.Range("H136").Value = 4174.154
$ws.Range("I136").Value = 3751.2727
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 11253.8181
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -8703.8181
$ws.Range("N136").Value = -24600
